# Update the "想去人数" (interested-count) column F values across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets to
# match the freshly scraped figures (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- 展览 -------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 1235
$ws1.Range("F11").Value = 285
$ws1.Range("F12").Value = 1084
$ws1.Range("F14").Value = 6848
$ws1.Range("F18").Value = 7756
$ws1.Range("F20").Value = 42
$ws1.Range("F21").Value = 4551
$ws1.Range("F23").Value = 2254
$ws1.Range("F25").Value = 4535
$ws1.Range("F26").Value = 240
$ws1.Range("F29").Value = 5
$ws1.Range("F30").Value = 265
$ws1.Range("F31").Value = 226
$ws1.Range("F32").Value = 5
$ws1.Range("F33").Value = 1946
$ws1.Range("F35").Value = 220
$ws1.Range("F37").Value = 524
$ws1.Range("F39").Value = 1340
$ws1.Range("F40").Value = 15
$ws1.Range("F41").Value = 2048
$ws1.Range("F42").Value = 2170

# --- 演出 -------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 10

# --- 全部类型 ---------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 1235
$ws4.Range("F13").Value = 285
$ws4.Range("F14").Value = 1084
$ws4.Range("F16").Value = 6848
$ws4.Range("F20").Value = 7756
$ws4.Range("F22").Value = 42
$ws4.Range("F23").Value = 4551
$ws4.Range("F25").Value = 2254
$ws4.Range("F27").Value = 4535
$ws4.Range("F28").Value = 240
$ws4.Range("F32").Value = 5
$ws4.Range("F34").Value = 265
$ws4.Range("F35").Value = 5
$ws4.Range("F36").Value = 1946
$ws4.Range("F38").Value = 220
$ws4.Range("F40").Value = 524
$ws4.Range("F42").Value = 10
$ws4.Range("F43").Value = 1340
$ws4.Range("F44").Value = 15
$ws4.Range("F45").Value = 2048
$ws4.Range("F47").Value = 2170
